# Auto-generated Excel COM-interop edit script
# Applies numeric value updates to specific cells across multiple sheets
# of the Moogle_Profits workbook, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 658.2
$ws.Range("I12").Value = 658.2
$ws.Range("K12").Value = 658.2
$ws.Range("M12").Value = -488.2

$ws.Range("H116").Value = 3961.5
$ws.Range("I116").Value = 3011
$ws.Range("J116").Value = 4912
$ws.Range("K116").Value = 3011
$ws.Range("L116").Value = 4912
$ws.Range("M116").Value = 431
$ws.Range("N116").Value = -11796

$ws.Range("H132").Value = 2615.984
$ws.Range("I132").Value = 2368.4385
$ws.Range("K132").Value = 7105.315500000001
$ws.Range("M132").Value = -4575.315500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 21742618
$ws.Range("I61").Value = 2963.238
$ws.Range("K61").Value = 2963.238
$ws.Range("M61").Value = -2751.238

$ws.Range("H63").Value = 3637.5386
$ws.Range("I63").Value = 3098.625
$ws.Range("J63").Value = 4499.8
$ws.Range("K63").Value = 3098.625
$ws.Range("L63").Value = 4499.8
$ws.Range("M63").Value = -2412.625
$ws.Range("N63").Value = -5871.8

$ws.Range("H66").Value = 3637.5386
$ws.Range("I66").Value = 3098.625
$ws.Range("J66").Value = 4499.8
$ws.Range("K66").Value = 15493.125
$ws.Range("L66").Value = 22499
$ws.Range("M66").Value = -12061.125
$ws.Range("N66").Value = -29363

$ws.Range("H74").Value = 5599.8237
$ws.Range("I74").Value = 3171.2144
$ws.Range("J74").Value = 16933.334
$ws.Range("K74").Value = 3171.2144
$ws.Range("L74").Value = 16933.334
$ws.Range("M74").Value = -2297.2144
$ws.Range("N74").Value = -18681.334

$ws.Range("H77").Value = 5599.8237
$ws.Range("I77").Value = 3171.2144
$ws.Range("J77").Value = 16933.334
$ws.Range("K77").Value = 15856.072
$ws.Range("L77").Value = 84666.67
$ws.Range("M77").Value = -11488.072
$ws.Range("N77").Value = -93402.67

$ws.Range("H136").Value = 21742618
$ws.Range("I136").Value = 2963.238
$ws.Range("K136").Value = 8889.714
$ws.Range("M136").Value = -6339.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 199950
$ws.Range("J76").Value = 199950
$ws.Range("L76").Value = 199950
$ws.Range("N76").Value = -200580

$ws.Range("H79").Value = 199950
$ws.Range("J79").Value = 199950
$ws.Range("L79").Value = 199950
$ws.Range("N79").Value = -202134

$ws.Range("H134").Value = 2214.5417
$ws.Range("I134").Value = 1499.1177
$ws.Range("K134").Value = 4497.3531
$ws.Range("M134").Value = -1962.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3386.6667
$ws.Range("J22").Value = 4449.6665
$ws.Range("L22").Value = 4449.6665
$ws.Range("N22").Value = -5149.6665

$ws.Range("H70").Value = 199950
$ws.Range("J70").Value = 199950
$ws.Range("L70").Value = 199950
$ws.Range("N70").Value = -200580

$ws.Range("H73").Value = 199950
$ws.Range("J73").Value = 199950
$ws.Range("L73").Value = 199950
$ws.Range("N73").Value = -202134

$ws.Range("H80").Value = 199950
$ws.Range("J80").Value = 199950
$ws.Range("L80").Value = 199950
$ws.Range("N80").Value = -202196

$ws.Range("H83").Value = 199950
$ws.Range("J83").Value = 199950
$ws.Range("L83").Value = 599850
$ws.Range("N83").Value = -611082

$ws.Range("H132").Value = 3989.975
$ws.Range("I132").Value = 3202.5925
$ws.Range("K132").Value = 9607.7775
$ws.Range("M132").Value = -7077.7775

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1822.6666
$ws.Range("I14").Value = 1822.6666
$ws.Range("K14").Value = 5467.9998
$ws.Range("M14").Value = -5294.9998

$ws.Range("H122").Value = 911.9167
$ws.Range("J122").Value = 896
$ws.Range("L122").Value = 8064
$ws.Range("N122").Value = -12964

$ws.Range("H131").Value = 791704.2
$ws.Range("I131").Value = 934.17645
$ws.Range("K131").Value = 2802.52935
$ws.Range("M131").Value = 2237.47065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 20439.6
$ws.Range("I7").Value = 24999
$ws.Range("J7").Value = 19299.75
$ws.Range("K7").Value = 24999
$ws.Range("L7").Value = 19299.75
$ws.Range("M7").Value = -24887
$ws.Range("N7").Value = -19523.75

$ws.Range("H8").Value = 20439.6
$ws.Range("I8").Value = 24999
$ws.Range("J8").Value = 19299.75
$ws.Range("K8").Value = 24999
$ws.Range("L8").Value = 19299.75
$ws.Range("M8").Value = -24860
$ws.Range("N8").Value = -19577.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4500.8
$ws.Range("I40").Value = 834.6667
$ws.Range("K40").Value = 834.6667
$ws.Range("M40").Value = -698.6667

$ws.Range("H68").Value = 7214.0586
$ws.Range("I68").Value = 8626.272000000001
$ws.Range("J68").Value = 4625
$ws.Range("K68").Value = 8626.272000000001
$ws.Range("L68").Value = 4625
$ws.Range("M68").Value = -7877.272000000001
$ws.Range("N68").Value = -6123

$ws.Range("H71").Value = 7214.0586
$ws.Range("I71").Value = 8626.272000000001
$ws.Range("J71").Value = 4625
$ws.Range("K71").Value = 43131.36
$ws.Range("L71").Value = 23125
$ws.Range("M71").Value = -39387.36
$ws.Range("N71").Value = -30613

$ws.Range("H122").Value = 5929.3125
$ws.Range("I122").Value = 6205.3335
$ws.Range("J122").Value = 5101.25
$ws.Range("K122").Value = 18616.0005
$ws.Range("L122").Value = 15303.75
$ws.Range("M122").Value = -16166.0005
$ws.Range("N122").Value = -20203.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 39999
$ws.Range("J63").Value = 46665.332
$ws.Range("L63").Value = 46665.332
$ws.Range("N63").Value = -47913.332

$ws.Range("H66").Value = 39999
$ws.Range("J66").Value = 46665.332
$ws.Range("L66").Value = 139995.996
$ws.Range("N66").Value = -146235.996

$ws.Range("H75").Value = 176625
$ws.Range("J75").Value = 199950
$ws.Range("L75").Value = 199950
$ws.Range("N75").Value = -201822

$ws.Range("H78").Value = 176625
$ws.Range("J78").Value = 199950
$ws.Range("L78").Value = 599850
$ws.Range("N78").Value = -609210

$ws.Range("H132").Value = 2502.1
$ws.Range("I132").Value = 1758.1351
$ws.Range("J132").Value = 4619.5386
$ws.Range("K132").Value = 5274.4053
$ws.Range("L132").Value = 13858.6158
$ws.Range("M132").Value = -2744.4053
$ws.Range("N132").Value = -18918.6158
